$d = $word.ActiveDocument

# The old paragraph text (Greek, with bold/italic "2018" lead-in, mixed-run formatting)
# that must be collapsed into a single plain run with the translated dates.
$oldText = "2018 Ημερομηνίες παρατήρησης για τον αστερισμό του Περσεύς: 30 Οκτωβρίου-8 Νοεμβρίου και 29 Νοεμβρίου-8 Δεκεμβρίου"
$newText = "Ημερομηνίες παρατήρησης για τον αστερισμό του Perseus: 16-25 Ιανουαρίου, 7-16 Νοεμβρίου, 6-15 Δεκεμβρίου"

# First, locate every paragraph whose text starts with the old content (there are
# a few duplicated occurrences of this notice throughout the document).
$matchIdx = New-Object System.Collections.ArrayList
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith($oldText)) {
        [void]$matchIdx.Add($i)
    }
}

# Then rewrite each matching paragraph: wipe all of its runs and replace them
# with a single, unformatted run holding the new (translated) text.
foreach ($idx in $matchIdx) {
    $p = $d.Paragraphs.Item($idx)
    $rng = $p.Range
    $bodyRng = $d.Range($rng.Start, $rng.End - 1)
    $bodyRng.Text = ""
    $insertPoint = $d.Range($p.Range.Start, $p.Range.Start)
    $insertPoint.InsertAfter($newText)
}
